$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("G4611")

$row = 30

$ws.Cells.Item($row, 1).Value = "Philippines"
$ws.Cells.Item($row, 2).Value = "philippines"
$ws.Cells.Item($row, 3).Value = "Bicol River Basin"
$ws.Cells.Item($row, 4).Value = "bicol"
$ws.Cells.Item($row, 5).Value = "Nabua"
$ws.Cells.Item($row, 6).Value = "G4611"
$ws.Cells.Item($row, 7).Value = "primary"

# forecast_date must stay a literal text string (matches the rest of column H),
# not get auto-converted into a date serial number by COM's type inference.
$dateCell = $ws.Cells.Item($row, 8)
$dateCell.NumberFormat = "@"
$dateCell.Value2 = "2025-10-29"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 13.37499999999993
$ws.Cells.Item($row, 11).Value = 123.3249999999996
$ws.Cells.Item($row, 12).Value = 5
$ws.Cells.Item($row, 13).Value = 483.4323679605675
$ws.Cells.Item($row, 14).Value = "LOW"
$ws.Cells.Item($row, 15).Value = 348.6773053168241
$ws.Cells.Item($row, 16).Value = 483.4323679605675
$ws.Cells.Item($row, 17).Value = 50
$ws.Cells.Item($row, 18).Value = 0
$ws.Cells.Item($row, 19).Value = 0
$ws.Cells.Item($row, 20).Value = 64.51953125
$ws.Cells.Item($row, 21).Value = 70.49047088623047
$ws.Cells.Item($row, 22).Value = 43.9609375
$ws.Cells.Item($row, 23).Value = 108.96875
$ws.Cells.Item($row, 24).Value = 58.53125
$ws.Cells.Item($row, 25).Value = 80.849609375
$ws.Cells.Item($row, 26).Value = $false
$ws.Cells.Item($row, 27).Value = -86.65386607806478
